# Auto-generated from the OOXML diff: update D (Price) and E (Volume 1h) columns
# for the crypto rows, plus the two row-swaps (B/C/D/E for rows 42/43 and 48/49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.734.60'
$ws.Range("E2").Value = '  +4.50%  '
$ws.Range("D3").Value = '1.875.77'
$ws.Range("E3").Value = '  +2.54%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '''338.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.26%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("D7").Value = '''0.4714'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.80%  '
$ws.Range("D8").Value = '''0.4028'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.38%  '
$ws.Range("D9").Value = '''47.71'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.80%  '
$ws.Range("D10").Value = '''0.08061'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.18%  '
$ws.Range("D11").Value = '''1.008'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.52%  '
$ws.Range("D12").Value = '''22.22'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.33%  '
$ws.Range("D13").Value = '''6.064'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.09%  '
$ws.Range("D14").Value = '1.868.48'
$ws.Range("E14").Value = '  +2.81%  '
$ws.Range("D15").Value = '''7.302'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.58%  '
$ws.Range("D16").Value = '''90.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.73%  '
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").Value = '''0.00001046'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.84%  '
$ws.Range("D19").Value = '''0.06619'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '''17.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.05%  '
$ws.Range("D21").Value = '''1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").Value = '28.751.89'
$ws.Range("E22").Value = '  +4.64%  '
$ws.Range("D23").Value = '''5.512'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.29%  '
$ws.Range("E24").Value = '  +2.16%  '
$ws.Range("D25").Value = '''2.263'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("D26").Value = '2.088.81'
$ws.Range("E26").Value = '  +2.69%  '
$ws.Range("D27").Value = '''160.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.56%  '
$ws.Range("E28").Value = '  +2.12%  '
$ws.Range("D29").Value = '''2.136'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.35%  '
$ws.Range("D30").Value = '''5.520'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.55%  '
$ws.Range("D31").Value = '''120.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.60%  '
$ws.Range("D32").Value = '''0.9884'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.10%  '
$ws.Range("D33").Value = '''0.09562'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.48%  '
$ws.Range("D34").Value = '''3.660'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.24%  '
$ws.Range("E35").Value = '  +4.34%  '
$ws.Range("D36").Value = '''5.387'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.60%  '
$ws.Range("D37").Value = '''0.06202'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.67%  '
$ws.Range("D38").Value = '''0.02280'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.36%  '
$ws.Range("D39").Value = '''8.499'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.87%  '
$ws.Range("D40").Value = '''1.188'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("D41").Value = '''0.5977'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.64%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.1894'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.20%  '
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = '''0.9996'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").Value = '''10.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.20%  '
$ws.Range("D45").Value = '''1.280'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = '''12.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("D47").Value = '''0.5585'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.28%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''1.965'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.95%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.07278'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.88%  '
$ws.Range("D50").Value = '''2.115'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +13.63%  '
$ws.Range("D51").Value = '''112.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.63%  '
